# Auto-generated Excel COM-interop script
# Updates market-data columns (H-N) in the Goblin_Profits workbook sheets
# to match the scheduled-runner snapshot described in the commit diff.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 213.5
$ws.Range("I4").Value = 213.5
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 213.5
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -99.5
$ws.Range("H17").Value = 2389.15
$ws.Range("J17").Value = 2510.7778
$ws.Range("L17").Value = 7532.3334
$ws.Range("N17").Value = -7868.3334
$ws.Range("H33").Value = 748
$ws.Range("J33").Value = 1750.2858
$ws.Range("L33").Value = 1750.2858
$ws.Range("N33").Value = -2208.2858
$ws.Range("H64").Value = 9555.421
$ws.Range("J64").Value = 9937.5
$ws.Range("L64").Value = 9937.5
$ws.Range("N64").Value = -10433.5
$ws.Range("H67").Value = 9555.421
$ws.Range("J67").Value = 9937.5
$ws.Range("L67").Value = 9937.5
$ws.Range("N67").Value = -11653.5
$ws.Range("H75").Value = 250039380
$ws.Range("J75").Value = 250039380
$ws.Range("L75").Value = 250039380
$ws.Range("N75").Value = -250041252
$ws.Range("H78").Value = 250039380
$ws.Range("J78").Value = 250039380
$ws.Range("L78").Value = 750118140
$ws.Range("N78").Value = -750127500
$ws.Range("H80").Value = 38462890
$ws.Range("J80").Value = 1770.2
$ws.Range("L80").Value = 5310.6
$ws.Range("N80").Value = -7306.6
$ws.Range("H83").Value = 38462890
$ws.Range("J83").Value = 1770.2
$ws.Range("L83").Value = 15931.8
$ws.Range("N83").Value = -25915.8
$ws.Range("H86").Value = 3161.611
$ws.Range("I86").Value = 3462.0833
$ws.Range("J86").Value = 2560.6667
$ws.Range("K86").Value = 3462.0833
$ws.Range("L86").Value = 2560.6667
$ws.Range("M86").Value = -2339.0833
$ws.Range("N86").Value = -4806.6667
$ws.Range("H89").Value = 3161.611
$ws.Range("I89").Value = 3462.0833
$ws.Range("J89").Value = 2560.6667
$ws.Range("K89").Value = 17310.4165
$ws.Range("L89").Value = 12803.3335
$ws.Range("M89").Value = -11694.4165
$ws.Range("N89").Value = -24035.3335
$ws.Range("H92").Value = 1431.4
$ws.Range("I92").Value = 1366
$ws.Range("K92").Value = 1366
$ws.Range("M92").Value = -118
$ws.Range("H96").Value = 1095.1666
$ws.Range("J96").Value = 3000
$ws.Range("L96").Value = 9000
$ws.Range("N96").Value = -11746
$ws.Range("H107").Value = 1360.65
$ws.Range("I107").Value = 1313.2778
$ws.Range("K107").Value = 1313.2778
$ws.Range("M107").Value = 606.7221999999999
$ws.Range("H129").Value = 1929.9
$ws.Range("I129").Value = 974
$ws.Range("J129").Value = 2567.1667
$ws.Range("K129").Value = 2922
$ws.Range("L129").Value = 7701.500100000001
$ws.Range("M129").Value = 2078
$ws.Range("N129").Value = -17701.5001
$ws.Range("H132").Value = 3127100
$ws.Range("J132").Value = 10002249
$ws.Range("L132").Value = 30006747
$ws.Range("N132").Value = -30011807
$ws.Range("H141").Value = 5498.913
$ws.Range("J141").Value = 12983
$ws.Range("L141").Value = 38949
$ws.Range("N141").Value = -49309

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1989
$ws.Range("I2").Value = 840.25
$ws.Range("K2").Value = 840.25
$ws.Range("M2").Value = -727.25
$ws.Range("H32").Value = 7837.4595
$ws.Range("I32").Value = 7837.4595
$ws.Range("K32").Value = 7837.4595
$ws.Range("M32").Value = -7550.4595
$ws.Range("H45").Value = 2092.7273
$ws.Range("I45").Value = 1480
$ws.Range("K45").Value = 1480
$ws.Range("M45").Value = -1103
$ws.Range("H97").Value = 1121.15
$ws.Range("I97").Value = 914.1539
$ws.Range("K97").Value = 914.1539
$ws.Range("M97").Value = -418.1539
$ws.Range("H116").Value = 1989
$ws.Range("I116").Value = 840.25
$ws.Range("K116").Value = 840.25
$ws.Range("M116").Value = 1453.75
$ws.Range("H132").Value = 1726.2106
$ws.Range("I132").Value = 1726.2106
$ws.Range("K132").Value = 5178.6318
$ws.Range("M132").Value = -2648.6318

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1989
$ws.Range("I3").Value = 840.25
$ws.Range("K3").Value = 840.25
$ws.Range("M3").Value = -726.25
$ws.Range("H86").Value = 20834422
$ws.Range("I86").Value = 1124.7222
$ws.Range("K86").Value = 1124.7222
$ws.Range("M86").Value = -1.72219999999993
$ws.Range("H89").Value = 20834422
$ws.Range("I89").Value = 1124.7222
$ws.Range("K89").Value = 5623.611
$ws.Range("M89").Value = -7.610999999999876
$ws.Range("H107").Value = 5190.9287
$ws.Range("I107").Value = 3669.2222
$ws.Range("K107").Value = 3669.2222
$ws.Range("M107").Value = -1749.2222
$ws.Range("H134").Value = 626641
$ws.Range("I134").Value = 1573.3055
$ws.Range("J134").Value = 6252250
$ws.Range("K134").Value = 4719.916499999999
$ws.Range("L134").Value = 18756750
$ws.Range("M134").Value = -2184.916499999999
$ws.Range("N134").Value = -18761820

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 208.41667
$ws.Range("I7").Value = 208.41667
$ws.Range("K7").Value = 208.41667
$ws.Range("M7").Value = -95.41667000000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 34361690
$ws.Range("I4").Value = 40736140
$ws.Range("J4").Value = 11181866
$ws.Range("K4").Value = 122208420
$ws.Range("L4").Value = 33545598
$ws.Range("M4").Value = -122208308
$ws.Range("N4").Value = -33545822
$ws.Range("H5").Value = 1419.5883
$ws.Range("I5").Value = 1149.8462
$ws.Range("J5").Value = 2296.25
$ws.Range("K5").Value = 3449.5386
$ws.Range("L5").Value = 6888.75
$ws.Range("M5").Value = -3337.5386
$ws.Range("N5").Value = -7112.75
$ws.Range("H34").Value = 2333.2222
$ws.Range("I34").Value = 2599.8
$ws.Range("K34").Value = 7799.400000000001
$ws.Range("M34").Value = -7715.400000000001
$ws.Range("H39").Value = 3832.6667
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H55").Value = 373029.44
$ws.Range("I55").Value = 555768
$ws.Range("J55").Value = 44100
$ws.Range("K55").Value = 1667304
$ws.Range("L55").Value = 132300
$ws.Range("M55").Value = -1667127
$ws.Range("N55").Value = -132654
$ws.Range("H131").Value = 3177653.5
$ws.Range("I131").Value = 798.6667
$ws.Range("K131").Value = 2396.0001
$ws.Range("M131").Value = 2643.9999
$ws.Range("H135").Value = 1419.5883
$ws.Range("I135").Value = 1149.8462
$ws.Range("J135").Value = 2296.25
$ws.Range("K135").Value = 10348.6158
$ws.Range("L135").Value = 20666.25
$ws.Range("M135").Value = -7813.6158
$ws.Range("N135").Value = -25736.25

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 2000
$ws.Range("J5").Value = 2000
$ws.Range("L5").Value = 2000
$ws.Range("N5").Value = -2224
$ws.Range("H97").Value = 1343.2
$ws.Range("J97").Value = 5000
$ws.Range("L97").Value = 5000
$ws.Range("N97").Value = -5992
$ws.Range("H102").Value = 1595.4667
$ws.Range("I102").Value = 1138
$ws.Range("K102").Value = 1138
$ws.Range("M102").Value = 484
$ws.Range("H136").Value = 48496.535
$ws.Range("J136").Value = 48496.535
$ws.Range("L136").Value = 145489.605
$ws.Range("N136").Value = -150589.605

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5241.1577
$ws.Range("I7").Value = 4760.231
$ws.Range("J7").Value = 6283.1665
$ws.Range("K7").Value = 4760.231
$ws.Range("L7").Value = 6283.1665
$ws.Range("M7").Value = -4648.231
$ws.Range("N7").Value = -6507.1665
$ws.Range("H29").Value = 20016
$ws.Range("I29").Value = 20016
$ws.Range("K29").Value = 20016
$ws.Range("M29").Value = -19721
$ws.Range("H34").Value = 8374.75
$ws.Range("I34").Value = 8374.75
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 8374.75
$ws.Range("L34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -8202.75
$ws.Range("H43").Value = 23749.5
$ws.Range("J43").Value = 23749.5
$ws.Range("L43").Value = 23749.5
$ws.Range("N43").Value = -24135.5
$ws.Range("H100").Value = 34486384
$ws.Range("I100").Value = 27779784
$ws.Range("K100").Value = 27779784
$ws.Range("M100").Value = -27779243
$ws.Range("H126").Value = 5241.1577
$ws.Range("I126").Value = 4760.231
$ws.Range("J126").Value = 6283.1665
$ws.Range("K126").Value = 14280.693
$ws.Range("L126").Value = 18849.4995
$ws.Range("M126").Value = -11810.693
$ws.Range("N126").Value = -23789.4995

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 6683.3335
$ws.Range("I4").Value = 6683.3335
$ws.Range("K4").Value = 6683.3335
$ws.Range("M4").Value = -6570.3335
$ws.Range("H46").Value = 55499.5
$ws.Range("J46").Value = 55499.5
$ws.Range("L46").Value = 55499.5
$ws.Range("N46").Value = -55961.5
$ws.Range("H132").Value = 9012217
$ws.Range("I132").Value = 11113234
$ws.Range("K132").Value = 33339702
$ws.Range("M132").Value = -33337172
$ws.Range("H134").Value = 55499.5
$ws.Range("J134").Value = 55499.5
$ws.Range("L134").Value = 166498.5
$ws.Range("N134").Value = -171568.5
$ws.Range("H136").Value = 4781.5557
$ws.Range("I136").Value = 2783.1667
$ws.Range("K136").Value = 8349.5001000000001
$ws.Range("M136").Value = -5799.5001000000001

